$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column P (year 2022) added after existing column O (year 2021).
# For each data row, copy the formatting from the corresponding column O
# cell (so the new cells pick up the same existing style index) and then
# set the new value for 2022.

$newValues = @{
    4  = 2022
    5  = 11.4
    6  = 12.6
    7  = 9.8
    8  = 11.4
    9  = 5.4
    10 = 4.7
    11 = 3.4
    12 = 17.7
    13 = 20.5
    14 = 8.4
    16 = 12.9
    17 = 10.5
}

foreach ($row in 4,5,6,7,8,9,10,11,12,13,14,16,17) {
    $src = $ws.Range("O$row")
    $dst = $ws.Range("P$row")
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $dst.Value = $newValues[$row]
}

# Move the active selection to Q4 to match the saved view state.
$ws.Range("Q4").Select()
